$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the SUM_Shape_Length column (C) entirely - data + header
$ws.Range("C1:C10").EntireColumn.Delete()

# Update header text/labels
$ws.Range("A1").Value = "عرض معبر"
$ws.Range("B1").Value = "تعداد"

# Headers are no longer bold
$ws.Range("A1:B1").Font.Bold = $false
